# Trade #117 closed at 2026-02-16 21:44:30 - leadlag UP +0.000%
#
# Applies:
#  - Summary sheet: refreshed aggregate stats (Total Trades / Win Rate / P&L)
#  - leadlag sheet: rows 70-74 flip from OPEN -> CLOSED (time_exit_5min) and
#    a brand new OPEN trade (#117) is appended as row 92
#  - momentum sheet: row 26 flips from OPEN -> CLOSED (time_exit_5min)
#  - All Trades sheet: the 6 newly-closed trades are appended as rows 94-99
#  - Comparison sheet: refreshed per-strategy stats for leadlag / momentum

$wb = $excel.ActiveWorkbook

# Helper: write a value that must be stored as literal text even though it
# looks like a number/date/percentage (Excel's COM layer otherwise "smart"
# converts "72.4%" -> 0.724, "2026-02-16" -> a date, etc.)
function Set-Text($ws, $addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
}

function Set-Num($ws, $addr, $val) {
    $ws.Range($addr).Value = $val
}

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$sum = $wb.Worksheets.Item("Summary")

Set-Num  $sum "C2" 98
Set-Text $sum "D2" "72.4%"
Set-Text $sum "E2" "+31.0572%"
Set-Text $sum "F2" "+0.3169%"

Set-Num  $sum "C3" 90
Set-Text $sum "D3" "53.3%"
Set-Text $sum "E3" "+17.2466%"
Set-Text $sum "F3" "+0.1916%"

Set-Text $sum "D4" "88.5%"
Set-Text $sum "E4" "+13.8106%"
Set-Text $sum "F4" "+0.5312%"

# ---------------------------------------------------------------------------
# leadlag sheet - close out trades 93-97 (rows 70-74)
# ---------------------------------------------------------------------------
$ll = $wb.Worksheets.Item("leadlag")

Set-Num  $ll "G70" 67836.49142999999
Set-Text $ll "H70" "CLOSED"
Set-Num  $ll "I70" 0.8058
Set-Num  $ll "J70" 8.06
Set-Text $ll "M70" "time_exit_5min"
Set-Num  $ll "N70" 5

Set-Num  $ll "G71" 68109.178749
Set-Text $ll "H71" "CLOSED"
Set-Num  $ll "I71" 0.331
Set-Num  $ll "J71" 3.31
Set-Text $ll "M71" "time_exit_5min"
Set-Num  $ll "N71" 5

Set-Num  $ll "G72" 68522.985287
Set-Text $ll "H72" "CLOSED"
Set-Num  $ll "I72" -0.3265
Set-Num  $ll "J72" -3.26
Set-Text $ll "M72" "time_exit_5min"
Set-Num  $ll "N72" 5

Set-Num  $ll "G73" 68882.266298
Set-Text $ll "H73" "CLOSED"
Set-Num  $ll "I73" 0.7895
Set-Num  $ll "J73" 7.9
Set-Text $ll "M73" "time_exit_5min"
Set-Num  $ll "N73" 5

Set-Num  $ll "G74" 68529.14739500001
Set-Text $ll "H74" "CLOSED"
Set-Num  $ll "I74" 0.1783
Set-Num  $ll "J74" 1.78
Set-Text $ll "M74" "time_exit_5min"
Set-Num  $ll "N74" 5

# New trade #117 - appended as row 92 (still OPEN)
Set-Num  $ll "A92" 117
Set-Text $ll "B92" "2026-02-16"
Set-Text $ll "C92" "21:44:30"
Set-Text $ll "D92" "leadlag"
Set-Text $ll "E92" "UP"
Set-Num  $ll "F92" 68447.97500000001
Set-Text $ll "H92" "OPEN"
Set-Num  $ll "I92" 0
Set-Num  $ll "J92" 0
Set-Num  $ll "K92" 0.75
Set-Text $ll "L92" "Coinbase leading with 0.091% move"
Set-Num  $ll "N92" 0

# ---------------------------------------------------------------------------
# momentum sheet - close out trade 96 (row 26)
# ---------------------------------------------------------------------------
$mo = $wb.Worksheets.Item("momentum")

Set-Num  $mo "G26" 67530.375915
Set-Text $mo "H26" "CLOSED"
Set-Num  $mo "I26" 1.1252
Set-Num  $mo "J26" 11.25
Set-Text $mo "M26" "time_exit_5min"
Set-Num  $mo "N26" 5

# ---------------------------------------------------------------------------
# All Trades sheet - append the 6 trades that just closed (rows 94-99)
# ---------------------------------------------------------------------------
$all = $wb.Worksheets.Item("All Trades")

Set-Num  $all "A94" 93
Set-Text $all "B94" "2026-02-16"
Set-Text $all "C94" "21:38:43"
Set-Text $all "D94" "leadlag"
Set-Text $all "E94" "DOWN"
Set-Num  $all "F94" 68387.55
Set-Num  $all "G94" 67836.49142999999
Set-Text $all "H94" "CLOSED"
Set-Num  $all "I94" 0.8058
Set-Num  $all "J94" 8.06
Set-Num  $all "K94" 0.621
Set-Text $all "L94" "Binance leading with -0.062% move"
Set-Text $all "M94" "time_exit_5min"
Set-Num  $all "N94" 5

Set-Num  $all "A95" 94
Set-Text $all "B95" "2026-02-16"
Set-Text $all "C95" "21:38:54"
Set-Text $all "D95" "leadlag"
Set-Text $all "E95" "DOWN"
Set-Num  $all "F95" 68335.345
Set-Num  $all "G95" 68109.178749
Set-Text $all "H95" "CLOSED"
Set-Num  $all "I95" 0.331
Set-Num  $all "J95" 3.31
Set-Num  $all "K95" 0.75
Set-Text $all "L95" "Coinbase leading with -0.087% move"
Set-Text $all "M95" "time_exit_5min"
Set-Num  $all "N95" 5

Set-Num  $all "A96" 95
Set-Text $all "B96" "2026-02-16"
Set-Text $all "C96" "21:39:01"
Set-Text $all "D96" "leadlag"
Set-Text $all "E96" "DOWN"
Set-Num  $all "F96" 68300
Set-Num  $all "G96" 68522.985287
Set-Text $all "H96" "CLOSED"
Set-Num  $all "I96" -0.3265
Set-Num  $all "J96" -3.26
Set-Num  $all "K96" 0.75
Set-Text $all "L96" "Binance leading with -0.146% move"
Set-Text $all "M96" "time_exit_5min"
Set-Num  $all "N96" 5

Set-Num  $all "A97" 96
Set-Text $all "B97" "2026-02-16"
Set-Text $all "C97" "21:39:07"
Set-Text $all "D97" "momentum"
Set-Text $all "E97" "DOWN"
Set-Num  $all "F97" 68298.875
Set-Num  $all "G97" 67530.375915
Set-Text $all "H97" "CLOSED"
Set-Num  $all "I97" 1.1252
Set-Num  $all "J97" 11.25
Set-Num  $all "K97" 0.9
Set-Text $all "L97" "Downward momentum: -0.242% over 10 samples"
Set-Text $all "M97" "time_exit_5min"
Set-Num  $all "N97" 5

Set-Num  $all "A98" 97
Set-Text $all "B98" "2026-02-16"
Set-Text $all "C98" "21:39:14"
Set-Text $all "D98" "leadlag"
Set-Text $all "E98" "UP"
Set-Num  $all "F98" 68342.69500000001
Set-Num  $all "G98" 68882.266298
Set-Text $all "H98" "CLOSED"
Set-Num  $all "I98" 0.7895
Set-Num  $all "J98" 7.9
Set-Num  $all "K98" 0.75
Set-Text $all "L98" "Binance leading with 0.086% move"
Set-Text $all "M98" "time_exit_5min"
Set-Num  $all "N98" 5

Set-Num  $all "A99" 98
Set-Text $all "B99" "2026-02-16"
Set-Text $all "C99" "21:39:20"
Set-Text $all "D99" "leadlag"
Set-Text $all "E99" "UP"
Set-Num  $all "F99" 68407.14999999999
Set-Num  $all "G99" 68529.14739500001
Set-Text $all "H99" "CLOSED"
Set-Num  $all "I99" 0.1783
Set-Num  $all "J99" 1.78
Set-Num  $all "K99" 0.75
Set-Text $all "L99" "Binance leading with 0.172% move"
Set-Text $all "M99" "time_exit_5min"
Set-Num  $all "N99" 5

# ---------------------------------------------------------------------------
# Comparison sheet
# ---------------------------------------------------------------------------
$cmp = $wb.Worksheets.Item("Comparison")

Set-Num  $cmp "B2" 90
Set-Text $cmp "C2" "53.3%"
Set-Text $cmp "D2" "3.25"
Set-Text $cmp "E2" "+0.5187%"
Set-Text $cmp "F2" "-0.3060%"
Set-Text $cmp "G2" "1.69"

Set-Text $cmp "C3" "88.5%"
Set-Text $cmp "D3" "13.28"
Set-Text $cmp "E3" "+0.6493%"
Set-Text $cmp "G3" "1.16"
